# A new weekly price record was added to the series. It is inserted as
# the new row 4 (right after the two already-unchanged rows), pushing
# every following record down by one row (old row 4 -> new row 5, ...,
# old row 85 -> new row 86).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 4. Excel shifts rows
# 4:85 down to 5:86 and copies the row-4 formatting (incl. the date
# number format on column D) onto the freshly inserted row.
$ws.Rows.Item(4).Insert()

# The categorical columns (market/region/product taxonomy + unit +
# origin + kg-per-unit) are identical for every record in this sheet,
# so copy them straight from the row directly below (the record that
# used to be row 4, now row 5).
$ws.Range("A4").Value2 = $ws.Range("A5").Value2
$ws.Range("B4").Value2 = $ws.Range("B5").Value2
$ws.Range("C4").Value2 = $ws.Range("C5").Value2
$ws.Range("E4").Value2 = $ws.Range("E5").Value2
$ws.Range("F4").Value2 = $ws.Range("F5").Value2
$ws.Range("G4").Value2 = $ws.Range("G5").Value2
$ws.Range("H4").Value2 = $ws.Range("H5").Value2
$ws.Range("I4").Value2 = $ws.Range("I5").Value2
$ws.Range("J4").Value2 = $ws.Range("J5").Value2
$ws.Range("K4").Value2 = $ws.Range("K5").Value2
$ws.Range("L4").Value2 = $ws.Range("L5").Value2
$ws.Range("Q4").Value2 = $ws.Range("Q5").Value2
$ws.Range("R4").Value2 = $ws.Range("R5").Value2
$ws.Range("T4").Value2 = $ws.Range("T5").Value2

# New record's own data.
$ws.Range("D4").Value2 = 44817
$ws.Range("M4").Value2 = 150
$ws.Range("N4").Value2 = 24000
$ws.Range("O4").Value2 = 25000
$ws.Range("P4").Value2 = 24500
$ws.Range("S4").Value2 = 1225
